# "Got the colored letters to work" -- replace the old "UpBrack"
# helper column (and its bracket-template strings) with new StimSet
# labels that embed the colored-letter placeholders directly
# (e.g. "ABXD{EF}"), entered with a leading apostrophe so Excel
# stores them as literal text (quotePrefix) instead of trying to
# evaluate the curly braces.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop column B ("UpBrack") by shifting C:K left into B:J.  Using
# Copy/PasteSpecial (rather than Columns("B").Delete) keeps the
# column width metadata for the sheet clean.
$ws.Range("C1:K4").Copy()
$ws.Range("B1").PasteSpecial(-4104)
$ws.Columns("K").ClearContents()
$excel.CutCopyMode = $false

# Re-enter the formulas that moved with the shift so the correct
# shared-formula groups are rebuilt (row 2 stand-alone, rows 3:4
# shared, matching the original layout's pattern).
$ws.Range("E2").Formula = "=C2+B2"
$ws.Range("E3").Formula = "=C3+B3"
$ws.Range("E4").Formula = "=C4+B4"

$ws.Range("H2").Formula = "=B2+C2+D2"
$ws.Range("H3:H4").Formula = "=B3+C3+D3"

$ws.Range("J2").Formula = "=D2+I2"
$ws.Range("J3:J4").Formula = "=D3+I3"

# New StimSet labels with the colored-letter markers baked in.
# The leading "'" forces a quote-prefixed text cell.
$ws.Range("A2").Value = "'ABXD{EF}"
$ws.Range("A3").Value = "'{ABYDEF}"
$ws.Range("A4").Value = "'A{BZDE}F"

# Move the selection off the table, like in the saved file.
$ws.Range("A5").Select()
